# Added short pulse filter to filter out high speed anomalies
# Extends the speedo calibration table (Tabelle1!D:E) with four more
# speed/period rows (150/160/170/180 mph) and moves the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D30").Value = 150
$ws.Range("E30").Value = 6000

$ws.Range("D31").Value = 160
$ws.Range("E31").Value = 5625

$ws.Range("D32").Value = 170
$ws.Range("E32").Value = 5294

$ws.Range("D33").Value = 180
$ws.Range("E33").Value = 5000

$ws.Range("C30").Select() | Out-Null
